$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range('A19').Value = 111928767
$ws.Range('Q19').Value = 540472

# Row 20
$ws.Range('A20').Value = 111928807
$ws.Range('B20').Value = 56430
$ws.Range('E20').Value = 100109
$ws.Range('F20').Value = 'Tretåig hackspett'
$ws.Range('G20').Value = 'Picoides tridactylus'
$ws.Range('H20').Value = '(Linnaeus, 1758)'
$ws.Range('J20').ClearContents()
$ws.Range('L20').Value = ''
$ws.Range('M20').Value = 'färska spår'
$ws.Range('Q20').Value = 540500
$ws.Range('R20').Value = 7247569
$ws.Range('AF20').ClearContents()

# Row 21
$ws.Range('A21').Value = 111928431
$ws.Range('B21').Value = 89571
$ws.Range('E21').Value = 5432
$ws.Range('F21').Value = 'Granticka'
$ws.Range('G21').Value = 'Porodaedalea chrysoloma'
$ws.Range('H21').Value = '(Fr.) Fiasson & Niemelä'
$ws.Range('J21').Value = ''
$ws.Range('L21').ClearContents()
$ws.Range('M21').ClearContents()
$ws.Range('Q21').Value = 540531
$ws.Range('R21').Value = 7247629
$ws.Range('AF21').Value = ''

# Row 22
$ws.Range('A22').Value = 111927812
$ws.Range('B22').Value = 85850
$ws.Range('E22').Value = 510
$ws.Range('F22').Value = 'Doftskinn'
$ws.Range('G22').Value = 'Cystostereum murrayi'
$ws.Range('H22').Value = '(Berk. & M.A. Curtis.) Pouzar'
$ws.Range('Q22').Value = 540626
$ws.Range('R22').Value = 7247582
$ws.Range('AJ22').Value = 'gran'
$ws.Range('AK22').Value = 'Picea abies'
$ws.Range('AL22').Value = 'Granlåga'
$ws.Range('AO22').Value = 'Picea abies # Granlåga'

# Row 23
$ws.Range('A23').Value = 111927932
$ws.Range('B23').Value = 90235
$ws.Range('D23').Value = 'LC'
$ws.Range('E23').Value = 3298
$ws.Range('F23').Value = 'Trådticka'
$ws.Range('G23').Value = 'Climacocystis borealis'
$ws.Range('H23').Value = '(Fr.) Kotl. & Pouzar'
$ws.Range('J23').Value = ''
$ws.Range('L23').ClearContents()
$ws.Range('M23').ClearContents()
$ws.Range('Q23').Value = 540603
$ws.Range('R23').Value = 7247579
$ws.Range('AF23').Value = ''

# Row 24
$ws.Range('A24').Value = 111929133
$ws.Range('B24').Value = 89571
$ws.Range('E24').Value = 5432
$ws.Range('F24').Value = 'Granticka'
$ws.Range('G24').Value = 'Porodaedalea chrysoloma'
$ws.Range('H24').Value = '(Fr.) Fiasson & Niemelä'
$ws.Range('Q24').Value = 540589
$ws.Range('R24').Value = 7247593
$ws.Range('AJ24').ClearContents()
$ws.Range('AK24').ClearContents()
$ws.Range('AL24').ClearContents()
$ws.Range('AO24').ClearContents()

# Row 25
$ws.Range('A25').Value = 111928641
$ws.Range('B25').Value = 77650
$ws.Range('E25').Value = 6425
$ws.Range('F25').Value = 'Garnlav'
$ws.Range('G25').Value = 'Alectoria sarmentosa'
$ws.Range('H25').Value = '(Ach.) Ach.'
$ws.Range('Q25').Value = 540501
$ws.Range('R25').Value = 7247613
$ws.Range('AC25').Value = 'Förekommer i området'

# Row 26
$ws.Range('B26').Value = 90235

# Row 27
$ws.Range('A27').Value = 111928864
$ws.Range('B27').Value = 89553
$ws.Range('D27').Value = 'NT'
$ws.Range('E27').Value = 1202
$ws.Range('F27').Value = 'Ullticka'
$ws.Range('G27').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H27').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range('Q27').Value = 540588
$ws.Range('R27').Value = 7247583

# Row 28
$ws.Range('A28').Value = 111928182
$ws.Range('B28').Value = 56430
$ws.Range('E28').Value = 100109
$ws.Range('F28').Value = 'Tretåig hackspett'
$ws.Range('G28').Value = 'Picoides tridactylus'
$ws.Range('H28').Value = '(Linnaeus, 1758)'
$ws.Range('J28').ClearContents()
$ws.Range('L28').Value = ''
$ws.Range('M28').Value = 'färska spår'
$ws.Range('Q28').Value = 540578
$ws.Range('R28').Value = 7247609
$ws.Range('AC28').ClearContents()
$ws.Range('AF28').ClearContents()

# Row 29
$ws.Range('B29').Value = 89571
